$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1755.2632
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1781.081
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 5343.242999999999
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -7559.242999999999

$ws.Range("H137").Value = 28104.342
$ws.Range("I137").Value = 1313.3334
$ws.Range("J137").Value = 93864.09
$ws.Range("K137").Value = 3940.0002
$ws.Range("L137").Value = 281592.27
$ws.Range("M137").Value = -1390.0002
$ws.Range("N137").Value = -286692.27

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16133966
$ws.Range("I32").Value = 17860482
$ws.Range("K32").Value = 17860482
$ws.Range("M32").Value = -17860195

$ws.Range("H61").Value = 1401.8857
$ws.Range("I61").Value = 918.5517
$ws.Range("J61").Value = 3738
$ws.Range("K61").Value = 918.5517
$ws.Range("L61").Value = 3738
$ws.Range("M61").Value = -706.5517
$ws.Range("N61").Value = -4162

$ws.Range("H74").Value = 2267.8667
$ws.Range("I74").Value = 2662.5625
$ws.Range("J74").Value = 1816.7858
$ws.Range("K74").Value = 2662.5625
$ws.Range("L74").Value = 1816.7858
$ws.Range("M74").Value = -1788.5625
$ws.Range("N74").Value = -3564.7858

$ws.Range("H77").Value = 2267.8667
$ws.Range("I77").Value = 2662.5625
$ws.Range("J77").Value = 1816.7858
$ws.Range("K77").Value = 13312.8125
$ws.Range("L77").Value = 9083.929
$ws.Range("M77").Value = -8944.8125
$ws.Range("N77").Value = -17819.929

$ws.Range("H136").Value = 1401.8857
$ws.Range("I136").Value = 918.5517
$ws.Range("J136").Value = 3738
$ws.Range("K136").Value = 2755.6551
$ws.Range("L136").Value = 11214
$ws.Range("M136").Value = -205.6550999999999
$ws.Range("N136").Value = -16314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1654.3549
$ws.Range("I86").Value = 1788.0555
$ws.Range("J86").Value = 1469.2307
$ws.Range("K86").Value = 1788.0555
$ws.Range("L86").Value = 1469.2307
$ws.Range("M86").Value = -665.0554999999999
$ws.Range("N86").Value = -3715.2307

$ws.Range("H89").Value = 1654.3549
$ws.Range("I89").Value = 1788.0555
$ws.Range("J89").Value = 1469.2307
$ws.Range("K89").Value = 8940.2775
$ws.Range("L89").Value = 7346.1535
$ws.Range("M89").Value = -3324.2775
$ws.Range("N89").Value = -18578.1535

$ws.Range("H134").Value = 1534.0807
$ws.Range("I134").Value = 1424.3658
$ws.Range("J134").Value = 1748.2858
$ws.Range("K134").Value = 4273.097400000001
$ws.Range("L134").Value = 5244.857400000001
$ws.Range("M134").Value = -1738.097400000001
$ws.Range("N134").Value = -10314.8574

$ws.Range("H137").Value = 48000
$ws.Range("J137").Value = 48000
$ws.Range("L137").Value = 48000
$ws.Range("N137").Value = -58200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2270.2917
$ws.Range("I31").Value = 1377.2972
$ws.Range("J31").Value = 5274
$ws.Range("K31").Value = 1377.2972
$ws.Range("L31").Value = 5274
$ws.Range("M31").Value = -1082.2972
$ws.Range("N31").Value = -5864

$ws.Range("H34").Value = 2270.2917
$ws.Range("I34").Value = 1377.2972
$ws.Range("J34").Value = 5274
$ws.Range("K34").Value = 1377.2972
$ws.Range("L34").Value = 5274
$ws.Range("M34").Value = -1175.2972
$ws.Range("N34").Value = -5678

$ws.Range("H58").Value = 1611.2325
$ws.Range("I58").Value = 903.9706
$ws.Range("J58").Value = 4283.1113
$ws.Range("K58").Value = 903.9706
$ws.Range("L58").Value = 4283.1113
$ws.Range("M58").Value = -700.9706
$ws.Range("N58").Value = -4689.1113

$ws.Range("H132").Value = 1695.8914
$ws.Range("I132").Value = 1285.2188
$ws.Range("J132").Value = 2634.5715
$ws.Range("K132").Value = 3855.6564
$ws.Range("L132").Value = 7903.7145
$ws.Range("M132").Value = -1325.6564
$ws.Range("N132").Value = -12963.7145

$ws.Range("H134").Value = 1796.7609
$ws.Range("I134").Value = 1164.421
$ws.Range("J134").Value = 4800.375
$ws.Range("K134").Value = 3493.263
$ws.Range("L134").Value = 14401.125
$ws.Range("M134").Value = -958.2629999999999
$ws.Range("N134").Value = -19471.125

$ws.Range("H136").Value = 1611.2325
$ws.Range("I136").Value = 903.9706
$ws.Range("J136").Value = 4283.1113
$ws.Range("K136").Value = 2711.9118
$ws.Range("L136").Value = 12849.3339
$ws.Range("M136").Value = -161.9117999999999
$ws.Range("N136").Value = -17949.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3401.5386
$ws.Range("I3").Value = 2246.6667
$ws.Range("K3").Value = 6740.000100000001
$ws.Range("M3").Value = -6628.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4406.0625
$ws.Range("I70").Value = 4179.6
$ws.Range("K70").Value = 4179.6
$ws.Range("M70").Value = -3909.6

$ws.Range("H73").Value = 4406.0625
$ws.Range("I73").Value = 4179.6
$ws.Range("K73").Value = 4179.6
$ws.Range("M73").Value = -3243.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1464.36
$ws.Range("I132").Value = 1658.0632
$ws.Range("J132").Value = 735.6667
$ws.Range("K132").Value = 4974.1896
$ws.Range("L132").Value = 2207.0001
$ws.Range("M132").Value = -2444.1896
$ws.Range("N132").Value = -7267.0001

$ws.Range("H136").Value = 2101.2354
$ws.Range("I136").Value = 1578.2821
$ws.Range("J136").Value = 3800.8333
$ws.Range("K136").Value = 4734.846299999999
$ws.Range("L136").Value = 11402.4999
$ws.Range("M136").Value = -2184.846299999999
$ws.Range("N136").Value = -16502.4999

$ws.Range("H140").Value = 56000
$ws.Range("J140").Value = 56000
$ws.Range("L140").Value = 56000
$ws.Range("N140").Value = -66360

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 914.7143
$ws.Range("I113").Value = 914.7143
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2744.1429
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -574.1428999999998
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 1680.4445
$ws.Range("I132").Value = 1056.2693
$ws.Range("J132").Value = 2260.0356
$ws.Range("K132").Value = 3168.8079
$ws.Range("L132").Value = 6780.1068
$ws.Range("M132").Value = -638.8078999999998
$ws.Range("N132").Value = -11840.1068

$ws.Range("H136").Value = 2856.4822
$ws.Range("I136").Value = 2691.681
$ws.Range("J136").Value = 3717.111
$ws.Range("K136").Value = 8075.043
$ws.Range("L136").Value = 11151.333
$ws.Range("M136").Value = -5525.043
$ws.Range("N136").Value = -16251.333
